$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-6
# from 2023-09-16 (45185) to 2023-10-05 (45204)
$ws.Range("C2:C6").Value = 45204
